$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 110, shifting existing rows
# 110-128 down to 111-129 (formatting is carried down from the row above,
# matching Excel's native "Insert" behaviour).
$ws.Rows("110").Insert()

# Populate the newly inserted row 110 with the new price-report record.
$ws.Cells.Item(110, 1).Value = 5
$ws.Cells.Item(110, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(110, 3).Value = "Maule"
$ws.Cells.Item(110, 4).Value = 44511
$ws.Cells.Item(110, 5).Value = 7
$ws.Cells.Item(110, 6).Value = 100112021
$ws.Cells.Item(110, 7).Value = "Ají"
$ws.Cells.Item(110, 8).Value = "Americana (o)"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 200
$ws.Cells.Item(110, 11).Value = 25000
$ws.Cells.Item(110, 12).Value = 25000
$ws.Cells.Item(110, 13).Value = 25000
$ws.Cells.Item(110, 14).Value = "$/caja 14 kilos"
$ws.Cells.Item(110, 15).Value = "Región del Maule"
$ws.Cells.Item(110, 16).Value = 1786
$ws.Cells.Item(110, 17).Value = 14
$ws.Cells.Item(110, 18).Value = "Hortaliza"
